# Fixes to incomplete parameter files.
#
# Adds two missing parameter rows (R_auto, POC_lc) that were previously
# blank placeholder rows in the sheet, and renames the "ObservedMAR"
# parameter to "ObservedMAR_oc" (its value/units are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 was a blank row between "DOC_precip" (row 18) and the
# "#Min/Resp Parameters" section header (row 20). Fill it in with the new
# R_auto parameter.
$ws.Cells.Item(19, 1).Value2 = "R_auto"
$ws.Cells.Item(19, 2).Value2 = 0.8
$ws.Cells.Item(19, 3).Value2 = "unitless"

# Row 23 was a blank row between "RespParam" (row 22) and the "#Sediment"
# section header (row 24). Fill it in with the new POC_lc parameter.
$ws.Cells.Item(23, 1).Value2 = "POC_lc"
$ws.Cells.Item(23, 2).Value2 = 0.01
$ws.Cells.Item(23, 3).Value2 = "1/days"

# Rename the "ObservedMAR" parameter label to "ObservedMAR_oc" (value 249
# and units g/m2/yr stay the same).
$ws.Cells.Item(25, 1).Value2 = "ObservedMAR_oc"

# Match the author's final selection/view state.
[void]$ws.Range("G20").Select()
